$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" values under duplicate_image_filename (column E) for rows 2-21
$ws.Range("E2:E21").Value = "NA"
